$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$testName = "Fill in the required information completely"

$ws.Range("A12").Value = $testName
$ws.Range("B12").Value = "PASSED"
$ws.Range("C12").Value = "chrome"

$ws.Range("A13").Value = $testName
$ws.Range("B13").Value = "FAILED"
$ws.Range("C13").Value = "chrome"

$ws.Range("A14").Value = $testName
$ws.Range("B14").Value = "FAILED"
$ws.Range("C14").Value = "chrome"
